# Auto-generated edit script: updates market-price-derived columns (H-N)
# across ALC/BSM/CRP/CUL/GSM/LTW/WVR sheets per refreshed source data.
$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 115.833336
$ws.Range("I12").Value = 99
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 99
$ws.Range("L12").Value = 200
$ws.Range("M12").Value = 71
$ws.Range("N12").Value = -540
# Row 19
$ws.Range("H19").Value = 1020.23334
$ws.Range("I19").Value = 882.5294
$ws.Range("J19").Value = 1200.3077
$ws.Range("K19").Value = 882.5294
$ws.Range("L19").Value = 1200.3077
$ws.Range("M19").Value = -707.5294
$ws.Range("N19").Value = -1550.3077
# Row 55
$ws.Range("H55").Value = 326.35715
$ws.Range("I55").Value = 312.375
$ws.Range("J55").Value = 345
$ws.Range("K55").Value = 312.375
$ws.Range("L55").Value = 345
$ws.Range("M55").Value = -98.375
$ws.Range("N55").Value = -773
# Row 82
$ws.Range("H82").Value = 358.66666
$ws.Range("I82").Value = 358.66666
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1075.99998
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -669.9999800000001
$ws.Range("N82").ClearContents()
# Row 85
$ws.Range("H85").Value = 358.66666
$ws.Range("I85").Value = 358.66666
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1075.99998
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 328.0000199999999
$ws.Range("N85").ClearContents()
# Row 118
$ws.Range("H118").Value = 3566.5
$ws.Range("I118").Value = 1526.6666
$ws.Range("J118").Value = 4978.6924
$ws.Range("K118").Value = 4579.9998
$ws.Range("L118").Value = 14936.0772
$ws.Range("M118").Value = -2922.9998
$ws.Range("N118").Value = -18250.0772
# Row 141
$ws.Range("H141").Value = 3077.95
$ws.Range("I141").Value = 2590.875
$ws.Range("K141").Value = 7772.625
$ws.Range("M141").Value = -2592.625

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 50
$ws.Range("H50").Value = 47326.668
$ws.Range("J50").Value = 47326.668
$ws.Range("L50").Value = 47326.668
$ws.Range("N50").Value = -48474.668
# Row 94
$ws.Range("H94").Value = 1143.409
$ws.Range("I94").Value = 1000.46155
$ws.Range("K94").Value = 1000.46155
$ws.Range("M94").Value = -549.46155
# Row 105
$ws.Range("H105").Value = 4434.355
$ws.Range("I105").Value = 3223.2222
$ws.Range("J105").Value = 4929.8184
$ws.Range("K105").Value = 3223.2222
$ws.Range("L105").Value = 4929.8184
$ws.Range("M105").Value = -1476.2222
$ws.Range("N105").Value = -8423.8184
# Row 109
$ws.Range("H109").Value = 20684
$ws.Range("J109").Value = 20684
$ws.Range("L109").Value = 20684
$ws.Range("N109").Value = -23458
# Row 134
$ws.Range("H134").Value = 2876.2632
$ws.Range("I134").Value = 1856.44
$ws.Range("J134").Value = 4837.4614
$ws.Range("K134").Value = 5569.32
$ws.Range("L134").Value = 14512.3842
$ws.Range("M134").Value = -3034.32
$ws.Range("N134").Value = -19582.3842

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 105
$ws.Range("H105").Value = 2368.125
$ws.Range("I105").Value = 1584.2858
$ws.Range("J105").Value = 2977.7778
$ws.Range("K105").Value = 1584.2858
$ws.Range("L105").Value = 2977.7778
$ws.Range("M105").Value = 162.7141999999999
$ws.Range("N105").Value = -6471.7778
# Row 107
$ws.Range("H107").Value = 525
$ws.Range("I107").Value = 437.25
$ws.Range("J107").Value = 700.5
$ws.Range("K107").Value = 437.25
$ws.Range("L107").Value = 700.5
$ws.Range("M107").Value = 1482.75
$ws.Range("N107").Value = -4540.5
# Row 132
$ws.Range("H132").Value = 3316.577
$ws.Range("I132").Value = 2647.1333
$ws.Range("J132").Value = 4229.4546
$ws.Range("K132").Value = 7941.3999
$ws.Range("L132").Value = 12688.3638
$ws.Range("M132").Value = -5411.3999
$ws.Range("N132").Value = -17748.3638
# Row 134
$ws.Range("H134").Value = 884036.5600000001
$ws.Range("I134").Value = 2266.52
$ws.Range("K134").Value = 6799.559999999999
$ws.Range("M134").Value = -4264.559999999999

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 245.33333
$ws.Range("I7").Value = 300
$ws.Range("J7").Value = 218
$ws.Range("K7").Value = 900
$ws.Range("L7").Value = 654
$ws.Range("M7").Value = -788
$ws.Range("N7").Value = -878
# Row 38
$ws.Range("H38").Value = 146.125
$ws.Range("I38").Value = 248.14285
$ws.Range("J38").Value = 66.77778000000001
$ws.Range("K38").Value = 744.4285500000001
$ws.Range("L38").Value = 200.33334
$ws.Range("M38").Value = -397.4285500000001
$ws.Range("N38").Value = -894.33334
# Row 39
$ws.Range("H39").Value = 523.38464
$ws.Range("J39").Value = 523.38464
$ws.Range("L39").Value = 1570.15392
$ws.Range("N39").Value = -2158.15392
# Row 41
$ws.Range("H41").Value = 1100
$ws.Range("J41").Value = 2000
$ws.Range("L41").Value = 6000
$ws.Range("N41").Value = -6676
# Row 42
$ws.Range("H42").Value = 6999.4
$ws.Range("J42").Value = 6999.4
$ws.Range("L42").Value = 20998.2
$ws.Range("N42").Value = -22066.2
# Row 44
$ws.Range("H44").Value = 1330.1666
$ws.Range("I44").Value = 192.16667
$ws.Range("J44").Value = 1899.1666
$ws.Range("K44").Value = 576.50001
$ws.Range("L44").Value = 5697.4998
$ws.Range("M44").Value = -178.50001
$ws.Range("N44").Value = -6493.4998
# Row 46
$ws.Range("H46").Value = 2338
$ws.Range("I46").Value = 200
$ws.Range("J46").Value = 3050.6667
$ws.Range("K46").Value = 600
$ws.Range("L46").Value = 9152.000100000001
$ws.Range("M46").Value = -509
$ws.Range("N46").Value = -9334.000100000001
# Row 48
$ws.Range("H48").Value = 1083.5
$ws.Range("J48").Value = 1088.2354
$ws.Range("L48").Value = 3264.7062
$ws.Range("N48").Value = -3764.7062
# Row 49
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
# Row 68
$ws.Range("H68").Value = 1075.75
$ws.Range("J68").Value = 1481.2
$ws.Range("L68").Value = 4443.6
$ws.Range("N68").Value = -6065.6
# Row 71
$ws.Range("H71").Value = 1075.75
$ws.Range("J71").Value = 1481.2
$ws.Range("L71").Value = 13330.8
$ws.Range("N71").Value = -21442.8
# Row 113
$ws.Range("H113").Value = 508
$ws.Range("I113").Value = 498.5
$ws.Range("J113").Value = 513.06665
$ws.Range("K113").Value = 1495.5
$ws.Range("L113").Value = 1539.19995
$ws.Range("M113").Value = 674.5
$ws.Range("N113").Value = -5879.19995
# Row 118
$ws.Range("H118").Value = 1479.5294
$ws.Range("J118").Value = 990.2
$ws.Range("L118").Value = 2970.6
$ws.Range("N118").Value = -5456.6
# Row 131
$ws.Range("H131").Value = 842.3728599999999
$ws.Range("J131").Value = 878.7174
$ws.Range("L131").Value = 2636.1522
$ws.Range("N131").Value = -12716.1522

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1314.6154
$ws.Range("I97").Value = 1419.8889
$ws.Range("K97").Value = 1419.8889
$ws.Range("M97").Value = -923.8888999999999
# Row 132
$ws.Range("H132").Value = 4953.846
$ws.Range("I132").Value = 3895.5
$ws.Range("K132").Value = 11686.5
$ws.Range("M132").Value = -9156.5

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1199.8572
$ws.Range("I22").Value = 350.5
$ws.Range("J22").Value = 1341.4166
$ws.Range("K22").Value = 350.5
$ws.Range("L22").Value = 1341.4166
$ws.Range("M22").Value = -55.5
$ws.Range("N22").Value = -1931.4166
# Row 27
$ws.Range("H27").Value = 1199.8572
$ws.Range("I27").Value = 350.5
$ws.Range("J27").Value = 1341.4166
$ws.Range("K27").Value = 350.5
$ws.Range("L27").Value = 1341.4166
$ws.Range("M27").Value = -243.5
$ws.Range("N27").Value = -1555.4166
# Row 40
$ws.Range("H40").Value = 7789
$ws.Range("I40").Value = 7338.5557
$ws.Range("J40").Value = 8599.799999999999
$ws.Range("K40").Value = 7338.5557
$ws.Range("L40").Value = 8599.799999999999
$ws.Range("M40").Value = -7202.5557
$ws.Range("N40").Value = -8871.799999999999
# Row 46
$ws.Range("H46").Value = 1586.875
$ws.Range("J46").Value = 1698.5714
$ws.Range("L46").Value = 1698.5714
$ws.Range("N46").Value = -2074.5714
# Row 68
$ws.Range("H68").Value = 1928.2632
$ws.Range("I68").Value = 1764.75
$ws.Range("J68").Value = 2208.5715
$ws.Range("K68").Value = 1764.75
$ws.Range("L68").Value = 2208.5715
$ws.Range("M68").Value = -1015.75
$ws.Range("N68").Value = -3706.5715
# Row 71
$ws.Range("H71").Value = 1928.2632
$ws.Range("I71").Value = 1764.75
$ws.Range("J71").Value = 2208.5715
$ws.Range("K71").Value = 8823.75
$ws.Range("L71").Value = 11042.8575
$ws.Range("M71").Value = -5079.75
$ws.Range("N71").Value = -18530.8575
# Row 139
$ws.Range("H139").Value = 39932.777
$ws.Range("J139").Value = 39843.125
$ws.Range("L139").Value = 39843.125
$ws.Range("N139").Value = -50123.125

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 1309.9615
$ws.Range("I113").Value = 484.6
$ws.Range("J113").Value = 2435.4546
$ws.Range("K113").Value = 1453.8
$ws.Range("L113").Value = 7306.3638
$ws.Range("M113").Value = 716.1999999999998
$ws.Range("N113").Value = -11646.3638

